$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 73.57143000000001
$ws.Range("I4").Value = 73.57143000000001
$ws.Range("K4").Value = 73.57143000000001
$ws.Range("M4").Value = 40.42856999999999
$ws.Range("H15").Value = 1310.3823
$ws.Range("I15").Value = 1310.3823
$ws.Range("K15").Value = 3931.1469
$ws.Range("M15").Value = -3762.1469
$ws.Range("H33").Value = 190.23077
$ws.Range("I33").Value = 190.23077
$ws.Range("K33").Value = 190.23077
$ws.Range("M33").Value = 38.76922999999999
$ws.Range("H40").Value = 9666.666999999999
$ws.Range("J40").Value = 16334
$ws.Range("L40").Value = 16334
$ws.Range("N40").Value = -16684
$ws.Range("H86").Value = 4493.4116
$ws.Range("I86").Value = 3421.2727
$ws.Range("J86").Value = 6459
$ws.Range("K86").Value = 3421.2727
$ws.Range("L86").Value = 6459
$ws.Range("M86").Value = -2298.2727
$ws.Range("N86").Value = -8705
$ws.Range("H89").Value = 4493.4116
$ws.Range("I89").Value = 3421.2727
$ws.Range("J89").Value = 6459
$ws.Range("K89").Value = 17106.3635
$ws.Range("L89").Value = 32295
$ws.Range("M89").Value = -11490.3635
$ws.Range("N89").Value = -43527
$ws.Range("H112").Value = 1783.875
$ws.Range("J112").Value = 4429.3335
$ws.Range("L112").Value = 13288.0005
$ws.Range("N112").Value = -15504.0005
$ws.Range("H116").Value = 14500.556
$ws.Range("I116").Value = 12249.75
$ws.Range("J116").Value = 16301.2
$ws.Range("K116").Value = 12249.75
$ws.Range("L116").Value = 16301.2
$ws.Range("M116").Value = -8807.75
$ws.Range("N116").Value = -23185.2
$ws.Range("H132").Value = 2078.138
$ws.Range("I132").Value = 1850.64
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 5551.92
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3021.92
$ws.Range("N132").Value = -15560
$ws.Range("H137").Value = 2362.0435
$ws.Range("I137").Value = 1762.52
$ws.Range("J137").Value = 3075.762
$ws.Range("K137").Value = 5287.559999999999
$ws.Range("L137").Value = 9227.286
$ws.Range("M137").Value = -2737.559999999999
$ws.Range("N137").Value = -14327.286

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2228.35
$ws.Range("I61").Value = 1739.4722
$ws.Range("K61").Value = 1739.4722
$ws.Range("M61").Value = -1527.4722
$ws.Range("H132").Value = 6139.5557
$ws.Range("I132").Value = 1664.6666
$ws.Range("J132").Value = 8377
$ws.Range("K132").Value = 4993.9998
$ws.Range("L132").Value = 25131
$ws.Range("M132").Value = -2463.9998
$ws.Range("N132").Value = -30191
$ws.Range("H133").Value = 75259.75
$ws.Range("J133").Value = 75259.75
$ws.Range("L133").Value = 75259.75
$ws.Range("N133").Value = -80319.75
$ws.Range("H136").Value = 2228.35
$ws.Range("I136").Value = 1739.4722
$ws.Range("K136").Value = 5218.4166
$ws.Range("M136").Value = -2668.4166

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1540
$ws.Range("J86").Value = 1800
$ws.Range("L86").Value = 1800
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 1540
$ws.Range("J89").Value = 1800
$ws.Range("L89").Value = 9000
$ws.Range("N89").Value = -20232

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4268.108
$ws.Range("I58").Value = 3337.6897
$ws.Range("J58").Value = 7640.875
$ws.Range("K58").Value = 3337.6897
$ws.Range("L58").Value = 7640.875
$ws.Range("M58").Value = -3134.6897
$ws.Range("N58").Value = -8046.875
$ws.Range("H99").Value = 2750
$ws.Range("I99").Value = 2750
$ws.Range("K99").Value = 2750
$ws.Range("M99").Value = -1252
$ws.Range("H122").Value = 7492.154
$ws.Range("I122").Value = 2649.7144
$ws.Range("J122").Value = 13141.667
$ws.Range("K122").Value = 7949.1432
$ws.Range("L122").Value = 39425.001
$ws.Range("M122").Value = -5499.1432
$ws.Range("N122").Value = -44325.001
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780
$ws.Range("H132").Value = 4019.6365
$ws.Range("I132").Value = 3792.0732
$ws.Range("J132").Value = 7129.6665
$ws.Range("K132").Value = 11376.2196
$ws.Range("L132").Value = 21388.9995
$ws.Range("M132").Value = -8846.2196
$ws.Range("N132").Value = -26448.9995
$ws.Range("H136").Value = 4268.108
$ws.Range("I136").Value = 3337.6897
$ws.Range("J136").Value = 7640.875
$ws.Range("K136").Value = 10013.0691
$ws.Range("L136").Value = 22922.625
$ws.Range("M136").Value = -7463.069100000001
$ws.Range("N136").Value = -28022.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 864.2857
$ws.Range("I26").Value = 933.3333
$ws.Range("J26").Value = 450
$ws.Range("K26").Value = 2799.9999
$ws.Range("L26").Value = 1350
$ws.Range("M26").Value = -2511.9999
$ws.Range("N26").Value = -1926
$ws.Range("H32").Value = 12958.667
$ws.Range("J32").Value = 12958.667
$ws.Range("L32").Value = 38876.001
$ws.Range("N32").Value = -39442.001
$ws.Range("H37").Value = 166666.67
$ws.Range("J37").Value = 166666.67
$ws.Range("L37").Value = 500000.01
$ws.Range("N37").Value = -500224.01
$ws.Range("H114").Value = 1051.7
$ws.Range("I114").Value = 736.8570999999999
$ws.Range("K114").Value = 2210.5713
$ws.Range("M114").Value = 1043.4287
$ws.Range("H117").Value = 2601.0908
$ws.Range("I117").Value = 702.8570999999999
$ws.Range("J117").Value = 5923
$ws.Range("K117").Value = 2108.5713
$ws.Range("L117").Value = 17769
$ws.Range("M117").Value = 1333.4287
$ws.Range("N117").Value = -24653
$ws.Range("H121").Value = 3061.2727
$ws.Range("I121").Value = 1466.1666
$ws.Range("J121").Value = 4975.4
$ws.Range("K121").Value = 4398.4998
$ws.Range("L121").Value = 14926.2
$ws.Range("M121").Value = -3088.4998
$ws.Range("N121").Value = -17546.2
$ws.Range("H132").Value = 4362.75
$ws.Range("I132").Value = 3168.625
$ws.Range("K132").Value = 28517.625
$ws.Range("M132").Value = -25987.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10620.2
$ws.Range("I80").Value = 1595
$ws.Range("K80").Value = 1595
$ws.Range("M80").Value = -597
$ws.Range("H83").Value = 10620.2
$ws.Range("I83").Value = 1595
$ws.Range("K83").Value = 7975
$ws.Range("M83").Value = -2983
$ws.Range("H113").Value = 3630.875
$ws.Range("I113").Value = 1850
$ws.Range("J113").Value = 4699.4
$ws.Range("K113").Value = 1850
$ws.Range("L113").Value = 4699.4
$ws.Range("M113").Value = 320
$ws.Range("N113").Value = -9039.4
$ws.Range("H122").Value = 7820
$ws.Range("I122").Value = 4453.5
$ws.Range("J122").Value = 9503.25
$ws.Range("K122").Value = 13360.5
$ws.Range("L122").Value = 28509.75
$ws.Range("M122").Value = -10910.5
$ws.Range("N122").Value = -33409.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5906.56
$ws.Range("I7").Value = 2556.4119
$ws.Range("K7").Value = 2556.4119
$ws.Range("M7").Value = -2444.4119
$ws.Range("H33").Value = 15003.5
$ws.Range("I33").Value = 21007
$ws.Range("K33").Value = 21007
$ws.Range("M33").Value = -20717
$ws.Range("H40").Value = 5060.4443
$ws.Range("I40").Value = 3448.4285
$ws.Range("K40").Value = 3448.4285
$ws.Range("M40").Value = -3312.4285
$ws.Range("H46").Value = 4471.2856
$ws.Range("I46").Value = 3400
$ws.Range("J46").Value = 4899.8
$ws.Range("K46").Value = 3400
$ws.Range("L46").Value = 4899.8
$ws.Range("M46").Value = -3212
$ws.Range("N46").Value = -5275.8
$ws.Range("H122").Value = 5138.724
$ws.Range("I122").Value = 4714.385
$ws.Range("K122").Value = 14143.155
$ws.Range("M122").Value = -11693.155
$ws.Range("H126").Value = 5906.56
$ws.Range("I126").Value = 2556.4119
$ws.Range("K126").Value = 7669.2357
$ws.Range("M126").Value = -5199.2357
$ws.Range("H130").Value = 66725
$ws.Range("J130").Value = 66725
$ws.Range("L130").Value = 66725
$ws.Range("N130").Value = -76765
$ws.Range("H132").Value = 4944
$ws.Range("I132").Value = 3214.3
$ws.Range("J132").Value = 7415
$ws.Range("K132").Value = 9642.900000000001
$ws.Range("L132").Value = 22245
$ws.Range("M132").Value = -7112.900000000001
$ws.Range("N132").Value = -27305
$ws.Range("H133").Value = 52730.4
$ws.Range("J133").Value = 52730.4
$ws.Range("L133").Value = 52730.4
$ws.Range("N133").Value = -57790.4
$ws.Range("H136").Value = 3495.3542
$ws.Range("I136").Value = 1897.7222
$ws.Range("J136").Value = 8288.25
$ws.Range("K136").Value = 5693.1666
$ws.Range("L136").Value = 24864.75
$ws.Range("M136").Value = -3143.1666
$ws.Range("N136").Value = -29964.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 31029
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 31029
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 31029
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -31435
$ws.Range("H122").Value = 8401.241
$ws.Range("I122").Value = 2619.0667
$ws.Range("J122").Value = 14596.429
$ws.Range("K122").Value = 7857.2001
$ws.Range("L122").Value = 43789.287
$ws.Range("M122").Value = -5407.2001
$ws.Range("N122").Value = -48689.287
$ws.Range("H126").Value = 3332.8333
$ws.Range("I126").Value = 1499
$ws.Range("K126").Value = 4497
$ws.Range("M126").Value = -2027
$ws.Range("H132").Value = 3105.6
$ws.Range("J132").Value = 3100.8333
$ws.Range("L132").Value = 9302.499899999999
$ws.Range("N132").Value = -14362.4999
$ws.Range("H136").Value = 2569.12
$ws.Range("I136").Value = 1140.1305
$ws.Range("J136").Value = 19002.5
$ws.Range("K136").Value = 3420.3915
$ws.Range("L136").Value = 57007.5
$ws.Range("M136").Value = -870.3914999999997
$ws.Range("N136").Value = -62107.5

Write-Host "Applied all Lamia_Profits updates"